$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99; this shifts existing rows 99..241 down to 100..242,
# preserving their values/styles (old row 241 becomes new row 242).
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with its data (same shape as the other rows,
# with Fecha = 44571 and Volumen = 4000; the remaining columns repeat the values
# that were already present for this record type).
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44571
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100114014
$ws.Range("G99").Value = "Betarraga"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 500
$ws.Range("M99").Value = 500
$ws.Range("N99").Value = "$/paquete 5 unidades"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 100
$ws.Range("Q99").Value = 5
$ws.Range("R99").Value = "Hortaliza"
